# tri_results.xlsx - "Added new triangulation2 script"
#
# The triangulation results sheet stores every value as text (the original
# file was produced by a script writing inline strings, not by typing into
# Excel), so every cell below has to be written back as text too - including
# the ones that look like plain numbers (columns D/E). Excel's COM object
# model auto-converts a numeric-looking string assigned to .Value into a
# real number, so for those we use the classic leading-apostrophe trick to
# force text storage, then ClearFormats() to drop the resulting "number
# stored as text" (quote-prefix) cell format so no formatting/style is left
# behind - only the displayed text changes, exactly like the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    # Plain assignment would leave a string like "(0.59,0.3)" alone, but a
    # purely numeric string would silently become a Double. Forcing it
    # through the apostrophe prefix keeps every cell a text cell regardless
    # of what it looks like, matching the workbook's existing convention.
    $ws.Range($cellRef).Value = "'" + $text
    $ws.Range($cellRef).ClearFormats()
}

# row, "Actual Co-Ord", "Estimated Co-Ord", "X Percentage Error", "Y Percentage Error"
$rows = @(
    @{ r = 2;  B = $null;              C = "0.5989,0.30008";   D = $null;       E = "-0.0163"   },
    @{ r = 3;  B = "0.3082,0.27783";   C = "0.30913,0.27897";  D = "-0.11603";  E = "-0.2278"   },
    @{ r = 4;  B = "0.05246,0.01999";  C = "0.05195,0.01774";  D = "0.06305";   E = "0.45022"   },
    @{ r = 5;  B = "0.39949,0.19011";  C = "0.39838,0.1891";   D = "0.13959";   E = "0.20268"   },
    @{ r = 6;  B = "0.78175,0.31904";  C = "0.78045,0.31426";  D = "0.16254";   E = "0.95531"   },
    @{ r = 7;  B = "0.2981,0.47984";   C = "0.2977,0.47797";   D = "0.04988";   E = "0.37359"   },
    @{ r = 8;  B = "0.50199,0.23276";  C = "0.5032,0.23243";   D = "-0.15121";  E = "0.06633"   },
    @{ r = 9;  B = "0.52727,0.4496";   C = "0.52899,0.44822";  D = "-0.21424";  E = "0.27622"   },
    @{ r = 10; B = "0.72516,0.45098";  C = "0.72779,0.45464";  D = "-0.32872";  E = "-0.73229"  },
    @{ r = 11; B = "0.56015,0.14296";  C = "0.56019,0.14323";  D = "-0.0053";   E = "-0.05451"  }
)

foreach ($row in $rows) {
    $r = $row.r
    if ($row.B) { Set-TextValue "B$r" ("(" + $row.B + ")") }
    if ($row.C) { Set-TextValue "C$r" ("(" + $row.C + ")") }
    if ($row.D) { Set-TextValue "D$r" $row.D }
    if ($row.E) { Set-TextValue "E$r" $row.E }
}
